# Re-order the data rows (2..13) of the single worksheet according to the
# fixed permutation below, and refresh the recomputed `fg_mf_total` (col X)
# value that changed for every data row (std computation / sample
# aggregation fix).
#
# new row (2..13) -> old row it should now contain
$rowMap = @{
    2  = 3
    3  = 6
    4  = 9
    5  = 5
    6  = 4
    7  = 8
    8  = 12
    9  = 7
    10 = 13
    11 = 11
    12 = 2
    13 = 10
}

$newFgMfTotal = 1.000031877590054

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 13
$firstCol = 1   # A
$lastCol = 24   # X
$nRows = $lastDataRow - $firstDataRow + 1
$nCols = $lastCol - $firstCol + 1

$srcRange = $ws.Range($ws.Cells.Item($firstDataRow, $firstCol), $ws.Cells.Item($lastDataRow, $lastCol))
$original = $srcRange.Value2

$reordered = New-Object 'object[,]' $nRows, $nCols
for ($newRow = $firstDataRow; $newRow -le $lastDataRow; $newRow++) {
    $oldRow = $rowMap[$newRow]
    $srcIdx = $oldRow - $firstDataRow + 1   # 1-based row index into $original
    for ($col = 1; $col -le $nCols; $col++) {
        $reordered[$newRow - $firstDataRow, $col - 1] = $original[$srcIdx, $col]
    }
    # column X is the last column (index $nCols) -> refresh fg_mf_total
    $reordered[$newRow - $firstDataRow, $nCols - 1] = $newFgMfTotal
}

$srcRange.Value2 = $reordered
